# plotEIC methods for fGroupsSet
#
# Insert a new "getEICsForFGroups" entry right above "getFeatures" (pushing
# it, and everything below, down by one row), and mark the (now shifted)
# "plotEIC" row as done ("G" column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "getFeatures" currently lives on row 19 - insert a fresh row above it so
# every following row (and the sheet's used range) shifts down by one.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the getEICsForFGroups entry:
# implement (D) = X, ionize (F) = X, done (G) = X
$ws.Range("A19").Value = "getEICsForFGroups"
$ws.Range("D19").Value = "X"
$ws.Range("F19").Value = "X"
$ws.Range("G19").Value = "X"

# "plotEIC" (originally row 33, now row 34 after the insert) is now done too
$ws.Range("G34").Value = "X"

# Mirror the author's final selection in the sheet view
$ws.Range("G35").Select()
